# The workbook is already open; grab the active workbook/worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The header in G1 was renamed from "Subscription Valid Till" to
# "Subscription Till". This is the only actual content change in the
# edit; everything else in the source diff is incidental metadata
# churn produced by re-saving the file.
$ws.Range("G1").Value = "Subscription Till"

$wb.Save()
